$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in Wednesday's hours for the "Sponsor Meeting" row (Time sheet up to Sep 29.2021)
$ws.Range("D8").Value = 1

# Update the active cell/selection to match the saved view state
$ws.Range("D8").Select()
